$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Valor Mora" total (E11) and the worker/period counters (C13, F13)
$ws.Range("E11").Value2 = 871688
$ws.Range("C13").Value2 = 2
$ws.Range("F13").Value2 = 15

# 2) Insert two new data rows right after the current last data row (row 29),
#    copying formatting+values from row 28 so the new rows inherit the normal
#    (non-bottom-border) style, and the existing last row (old row 29) slides
#    down to keep its bottom-border "last row" styling.
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(29).Insert()
$ws.Range("B28:J28").Copy($ws.Range("B29:J29"))
$ws.Range("B28:J28").Copy($ws.Range("B30:J30"))

# 3) New row 30: same worker (PAOLA ESTER NISPERUZA ESTRADA), new period 2509
$ws.Range("E30").Value2 = "2509"
$ws.Range("F30").Value2 = 52900
$ws.Range("G30").Value2 = 1322500
$ws.Range("H30").Value2 = ""
$ws.Range("I30").Value2 = ""
$ws.Range("J30").Value2 = ""

# 4) Row 31 (the old last row, now shifted down) becomes the new worker
#    REYNA ISABEL MIRANDA JARABA, also for period 2509
$ws.Range("C31").Value2 = "1047431464"
$ws.Range("D31").Value2 = "REYNA ISABEL MIRANDA JARABA"
$ws.Range("E31").Value2 = "2509"
$ws.Range("F31").Value2 = 78188
$ws.Range("G31").Value2 = 1954700
$ws.Range("H31").Value2 = ""
$ws.Range("I31").Value2 = ""
$ws.Range("J31").Value2 = ""

Write-Host "done"
